# Apply the "new orleans xlsx" update:
#  1. Insert a new "State" column into hotel_info right after "Hotel_Name"
#     (before "City"), populated with "Louisiana" for the existing hotel row.
#  2. Reorder the worksheet tabs so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

$hotelSheet  = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# --- 1. Insert the new State column -----------------------------------
# hotel_info header row: STR | Hotel_Name | City | Zip | TA_ReviewURL | ...
# Column C currently holds "City" -> insert a new blank column there so
# "City" (and everything after it) shifts right, then fill the new column.
$hotelSheet.Columns.Item(3).Insert()
$hotelSheet.Cells.Item(1, 3).Value = "State"
$hotelSheet.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. Reorder the sheet tabs -----------------------------------------
# Move review_info in front of hotel_info so the tab order becomes
# review_info, hotel_info.
$reviewSheet.Move($hotelSheet)
